$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-4.07%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'37.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-6.48%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'-1.08%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07707"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-6.16%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.347"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.71%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.918"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-6.72%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'8.177"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.33%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.013"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-10.26%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.72%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1189"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-13.36%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1857"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-7.34%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08679"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-4.59%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03396"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-3.60%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09696"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.01%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001365"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.33%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006074"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.72%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.616"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.83%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3408"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.83%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'5.022"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.35%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1268"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-4.04%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.02105"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'5,164.77%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04333"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.29%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001212"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.46%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004212"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-12.14%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001354"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'4.02%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02192"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-5.47%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04880"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-6.24%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007571"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.42%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.57%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-5.27%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002066"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.66%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008529"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.41%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006532"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.34%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'2.10%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001301"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-23.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("E51").Style = "Normal"
